$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove erroneous leaked cells (naive forecaster bug)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update recalculated forecast values (precision fix)
$ws.Range("E4").Value = 2.64925711235009
$ws.Range("C6").Value = 2.533533936850585
$ws.Range("E7").Value = 2.828066716168043
$ws.Range("E8").Value = 2.332261646026246
$ws.Range("C9").Value = 1.891592186533786
$ws.Range("E9").Value = 2.544631191216373
$ws.Range("C10").Value = 1.21254482274098
$ws.Range("E10").Value = 1.839804681163337
$ws.Range("E11").Value = 1.312870290004309
$ws.Range("E12").Value = 0.6705904529405782
$ws.Range("C13").Value = 0.618025493879526
$ws.Range("E13").Value = 0.7749619016294229
$ws.Range("C14").Value = 0.4712609263772816
$ws.Range("E14").Value = 0.8520644823059031
$ws.Range("E15").Value = 0.7487574275251818
$ws.Range("C17").Value = 1.905862317202112
$ws.Range("C20").Value = 4.109890522944326
$ws.Range("E20").Value = 3.628019428949014
$ws.Range("C21").Value = 1.917627847674042
$ws.Range("E21").Value = 2.69471174461664
$ws.Range("C24").Value = 1.119562422009124
$ws.Range("C27").Value = 2.252616573494315
$ws.Range("E27").Value = 1.46507610487594
$ws.Range("C28").Value = 1.344920716048215
$ws.Range("E28").Value = 1.037735724446587
$ws.Range("E29").Value = 1.751794502139248
$ws.Range("E31").Value = 2.057677568601424
$ws.Range("C32").Value = 2.195375580740744
$ws.Range("E34").Value = 3.933586883651397
$ws.Range("C35").Value = 2.014919551176164
$ws.Range("C36").Value = 2.491319804758541
$ws.Range("E36").Value = 2.78106797904647
$ws.Range("C37").Value = 2.358700676763137
$ws.Range("E37").Value = 3.271886281175829
$ws.Range("C38").Value = 2.777797690741446
$ws.Range("E38").Value = 2.073300717643911
$ws.Range("E40").Value = 1.872042068954638
$ws.Range("C41").Value = 2.740261495864793
$ws.Range("E41").Value = 4.569471876550879
$ws.Range("E42").Value = 1.589741018019186
$ws.Range("E43").Value = 1.167752392835819
$ws.Range("C44").Value = -4.149799191324066
$ws.Range("E44").Value = -2.475922651815632
$ws.Range("C45").Value = -1.678602239427673
$ws.Range("E45").Value = -0.5705663367256197
$ws.Range("C46").Value = -1.432689847121826
$ws.Range("E47").Value = -1.116764638702061
$ws.Range("C50").Value = 2.033479419175155
$ws.Range("E50").Value = 1.562315774899048
$ws.Range("C52").Value = 2.221594549876427
$ws.Range("C53").Value = 2.502458807779662
